# "update pw - flow case"
# Sheet1!D2 held the placeholder value "demo" (password/flow-case test data);
# update it to "smartmed" and move the active selection to D2, matching the
# updated test fixture for the flow-case test.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "smartmed"

$ws.Activate()
$ws.Range("D2").Select()
